$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 21
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3"
$ws.Range("I2").Value = "Real Estate"
$ws.Range("J2").Value = "C Company"
$ws.Range("K2").Value = "Berlin, Germany"
$ws.Range("M2").Value = "`n[2025-09-24 22:41:26] The customer confirms it's a good time to talk and inquires about the AI's knowledge of them. The AI recalls the customer's name, association with C company in the real estate industry in Berlin, Germany, and their dust allergy concerns."
